$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.725194398296594
$ws.Cells.Item(2, 4).Value = 0.2626496464037587
$ws.Cells.Item(2, 5).Value = 0.1765809838104246
$ws.Cells.Item(2, 6).Value = 5.593921383993347
$ws.Cells.Item(2, 7).Value = 0.002670869184767723
$ws.Cells.Item(2, 10).Value = 0.1929031843028355
$ws.Cells.Item(2, 12).Value = 1.27378031124033
$ws.Cells.Item(2, 13).Value = 0.922060958821433
$ws.Cells.Item(2, 14).Value = 2.657762862678169

$ws.Cells.Item(3, 2).Value = 2.661214646614837
$ws.Cells.Item(3, 4).Value = 0.2351416345204314
$ws.Cells.Item(3, 5).Value = 0.1535691909987733
$ws.Cells.Item(3, 6).Value = 5.536426800347783
$ws.Cells.Item(3, 7).Value = 0.002680822602941113
$ws.Cells.Item(3, 10).Value = 0.1682863146042592
$ws.Cells.Item(3, 12).Value = 1.223503659134337
$ws.Cells.Item(3, 13).Value = 0.8941446311664976
$ws.Cells.Item(3, 14).Value = 2.687968677205923

$ws.Cells.Item(4, 2).Value = 2.623542432703175
$ws.Cells.Item(4, 4).Value = 0.2184256650930649
$ws.Cells.Item(4, 5).Value = 0.139486579031022
$ws.Cells.Item(4, 6).Value = 5.504775140787615
$ws.Cells.Item(4, 7).Value = 0.002687245852607324
$ws.Cells.Item(4, 10).Value = 0.1531624957255104
$ws.Cells.Item(4, 12).Value = 1.193482695804789
$ws.Cells.Item(4, 13).Value = 0.8775837343299173
$ws.Cells.Item(4, 14).Value = 2.707391697807047

$ws.Cells.Item(5, 2).Value = 2.60859374123504
$ws.Cells.Item(5, 4).Value = 0.2116556328359138
$ws.Cells.Item(5, 5).Value = 0.1337582946326563
$ws.Cells.Item(5, 6).Value = 5.492787120187955
$ws.Cells.Item(5, 7).Value = 0.002689942126511437
$ws.Cells.Item(5, 10).Value = 0.146996042840513
$ws.Cells.Item(5, 12).Value = 1.181459948492687
$ws.Cells.Item(5, 13).Value = 0.8709796652599664
$ws.Cells.Item(5, 14).Value = 2.715526737100724

$ws.Cells.Item(6, 2).Value = 2.60613580192927
$ws.Cells.Item(6, 4).Value = 0.2105339472708749
$ws.Cells.Item(6, 5).Value = 0.1328077149705535
$ws.Cells.Item(6, 6).Value = 5.490851278509325
$ws.Cells.Item(6, 7).Value = 0.002690394607435797
$ws.Cells.Item(6, 10).Value = 0.1459718719240186
$ws.Cells.Item(6, 12).Value = 1.179476255555841
$ws.Cells.Item(6, 13).Value = 0.8698917662719765
$ws.Cells.Item(6, 14).Value = 2.716890829933913

$ws.Cells.Item(7, 2).Value = 2.623339200381793
$ws.Cells.Item(7, 4).Value = 0.218334195075073
$ws.Cells.Item(7, 5).Value = 0.1394092845354322
$ws.Cells.Item(7, 6).Value = 5.50460979033295
$ws.Cells.Item(7, 7).Value = 0.002687281896082394
$ws.Cells.Item(7, 10).Value = 0.1530793479128647
$ws.Cells.Item(7, 12).Value = 1.193319701562586
$ws.Cells.Item(7, 13).Value = 0.877494085465564
$ws.Cells.Item(7, 14).Value = 2.707500519541451

$ws.Cells.Item(8, 2).Value = 2.70279831129136
$ws.Cells.Item(8, 4).Value = 0.2531275917764901
$ws.Cells.Item(8, 5).Value = 0.1686360526341844
$ws.Cells.Item(8, 6).Value = 5.573334994001272
$ws.Cells.Item(8, 7).Value = 0.002674236612726272
$ws.Cells.Item(8, 10).Value = 0.1844164601699276
$ws.Cells.Item(8, 12).Value = 1.256267163948536
$ws.Cells.Item(8, 13).Value = 0.9123143355314411
$ws.Cells.Item(8, 14).Value = 2.667995658515537

$ws.Cells.Item(9, 2).Value = 2.871519075932156
$ws.Cells.Item(9, 4).Value = 0.3228289623567662
$ws.Cells.Item(9, 5).Value = 0.226380411838619
$ws.Cells.Item(9, 6).Value = 5.73743059908503
$ws.Cells.Item(9, 7).Value = 0.002651113346906607
$ws.Cells.Item(9, 10).Value = 0.2458514335800714
$ws.Cells.Item(9, 12).Value = 1.386566106772989
$ws.Cells.Item(9, 13).Value = 0.9852564945380067
$ws.Cells.Item(9, 14).Value = 2.597498925529656

$ws.Cells.Item(10, 2).Value = 3.003522902463146
$ws.Cells.Item(10, 4).Value = 0.3750676394433015
$ws.Cells.Item(10, 5).Value = 0.2691560250445804
$ws.Cells.Item(10, 6).Value = 5.876400966100789
$ws.Cells.Item(10, 7).Value = 0.002635601543884324
$ws.Cells.Item(10, 10).Value = 0.2910574689401528
$ws.Cells.Item(10, 12).Value = 1.486662115138301
$ws.Cells.Item(10, 13).Value = 1.041778231776263
$ws.Cells.Item(10, 14).Value = 2.549977778724269

$ws.Cells.Item(11, 2).Value = 3.065363142275714
$ws.Cells.Item(11, 4).Value = 0.3990853617507355
$ws.Cells.Item(11, 5).Value = 0.2887117528234597
$ws.Cells.Item(11, 6).Value = 5.943739360258093
$ws.Cells.Item(11, 7).Value = 0.002628860782518423
$ws.Cells.Item(11, 10).Value = 0.3116560605106145
$ws.Cells.Item(11, 12).Value = 1.533188131352631
$ws.Cells.Item(11, 13).Value = 1.068148605032604
$ws.Cells.Item(11, 14).Value = 2.529292008637228

$ws.Cells.Item(12, 2).Value = 3.089041156945939
$ws.Cells.Item(12, 4).Value = 0.408219163721725
$ws.Cells.Item(12, 5).Value = 0.2961325122042382
$ws.Cells.Item(12, 6).Value = 5.969840833304403
$ws.Cells.Item(12, 7).Value = 0.002626353250225883
$ws.Cells.Item(12, 10).Value = 0.3194625506817204
$ws.Cells.Item(12, 12).Value = 1.550952254346612
$ws.Cells.Item(12, 13).Value = 1.078230663539813
$ws.Cells.Item(12, 14).Value = 2.521593353562658

$ws.Cells.Item(13, 2).Value = 3.083930041415897
$ws.Cells.Item(13, 4).Value = 0.4062502730561732
$ws.Cells.Item(13, 5).Value = 0.294533608394218
$ws.Cells.Item(13, 6).Value = 5.964192492313686
$ws.Cells.Item(13, 7).Value = 0.002626891294154244
$ws.Cells.Item(13, 10).Value = 0.3177809835963217
$ws.Cells.Item(13, 12).Value = 1.547119900717576
$ws.Cells.Item(13, 13).Value = 1.076055010207924
$ws.Cells.Item(13, 14).Value = 2.523245400681475

$ws.Cells.Item(14, 2).Value = 3.067305909786171
$ws.Cells.Item(14, 4).Value = 0.3998360140564614
$ws.Cells.Item(14, 5).Value = 0.289321945210375
$ws.Cells.Item(14, 6).Value = 5.945874626620935
$ws.Cells.Item(14, 7).Value = 0.002628653585729678
$ws.Cells.Item(14, 10).Value = 0.3122981720664768
$ws.Cells.Item(14, 12).Value = 1.534646657650228
$ws.Cells.Item(14, 13).Value = 1.068976125881207
$ws.Cells.Item(14, 14).Value = 2.52865593462699

$ws.Cells.Item(15, 2).Value = 3.057157156641949
$ws.Cells.Item(15, 4).Value = 0.3959122181244368
$ws.Cells.Item(15, 5).Value = 0.2861317055411945
$ws.Cells.Item(15, 6).Value = 5.934733072112067
$ws.Cells.Item(15, 7).Value = 0.002629738896364771
$ws.Cells.Item(15, 10).Value = 0.3089406461465387
$ws.Cells.Item(15, 12).Value = 1.527025508357497
$ws.Cells.Item(15, 13).Value = 1.064652674459452
$ws.Cells.Item(15, 14).Value = 2.531987589534239

$ws.Cells.Item(16, 2).Value = 2.999517783818419
$ws.Cells.Item(16, 4).Value = 0.3735033308250593
$ws.Cells.Item(16, 5).Value = 0.2678800889790125
$ws.Cells.Item(16, 6).Value = 5.872083949829772
$ws.Cells.Item(16, 7).Value = 0.002636048390807336
$ws.Cells.Item(16, 10).Value = 0.2897121033421399
$ws.Cells.Item(16, 12).Value = 1.483641752098038
$ws.Cells.Item(16, 13).Value = 1.040068240682317
$ws.Cells.Item(16, 14).Value = 2.55134844426199

$ws.Cells.Item(17, 2).Value = 2.964618658988172
$ws.Cells.Item(17, 4).Value = 0.3598228826080003
$ws.Cells.Item(17, 5).Value = 0.2567091817190743
$ws.Cells.Item(17, 6).Value = 5.834712617329757
$ws.Cells.Item(17, 7).Value = 0.002639999661440089
$ws.Cells.Item(17, 10).Value = 0.2779257019856516
$ws.Cells.Item(17, 12).Value = 1.457283435142187
$ws.Cells.Item(17, 13).Value = 1.025156158064107
$ws.Cells.Item(17, 14).Value = 2.563464788458141

$ws.Cells.Item(18, 2).Value = 2.944714005874005
$ws.Cells.Item(18, 4).Value = 0.3519780270237902
$ws.Cells.Item(18, 5).Value = 0.2502929988728937
$ws.Cells.Item(18, 6).Value = 5.813604837787864
$ws.Cells.Item(18, 7).Value = 0.002642302057364299
$ws.Cells.Item(18, 10).Value = 0.2711496065923882
$ws.Cells.Item(18, 12).Value = 1.442215993434672
$ws.Cells.Item(18, 13).Value = 1.01664099694662
$ws.Cells.Item(18, 14).Value = 2.570521461180554

$ws.Cells.Item(19, 2).Value = 2.938003468347574
$ws.Cells.Item(19, 4).Value = 0.349325908956132
$ws.Cells.Item(19, 5).Value = 0.2481220951111709
$ws.Cells.Item(19, 6).Value = 5.806524329086187
$ws.Cells.Item(19, 7).Value = 0.002643086725289329
$ws.Cells.Item(19, 10).Value = 0.2688558345966499
$ws.Cells.Item(19, 12).Value = 1.437130324736643
$ws.Cells.Item(19, 13).Value = 1.013768491784759
$ws.Cells.Item(19, 14).Value = 2.572925769054116

$ws.Cells.Item(20, 2).Value = 2.968316280873921
$ws.Cells.Item(20, 4).Value = 0.3612767126100209
$ws.Cells.Item(20, 5).Value = 0.25789739882498
$ws.Cells.Item(20, 6).Value = 5.83865071598575
$ws.Cells.Item(20, 7).Value = 0.002639575967441344
$ws.Cells.Item(20, 10).Value = 0.2791800508904601
$ws.Cells.Item(20, 12).Value = 1.46007965857342
$ws.Cells.Item(20, 13).Value = 1.026737160245332
$ws.Cells.Item(20, 14).Value = 2.562165904934055

$ws.Cells.Item(21, 2).Value = 3.07218172773895
$ws.Cells.Item(21, 4).Value = 0.4017189656639175
$ws.Cells.Item(21, 5).Value = 0.2908523061363582
$ws.Cells.Item(21, 6).Value = 5.9512386164524
$ws.Cells.Item(21, 7).Value = 0.002628134738013717
$ws.Cells.Item(21, 10).Value = 0.313908426207604
$ws.Cells.Item(21, 12).Value = 1.538306371174656
$ws.Cells.Item(21, 13).Value = 1.071052743552741
$ws.Cells.Item(21, 14).Value = 2.527063071845369

$ws.Cells.Item(22, 2).Value = 3.141583377876032
$ws.Cells.Item(22, 4).Value = 0.4283776369121597
$ws.Cells.Item(22, 5).Value = 0.3124809022519059
$ws.Cells.Item(22, 6).Value = 6.028334119572605
$ws.Cells.Item(22, 7).Value = 0.002620919645905147
$ws.Cells.Item(22, 10).Value = 0.3366425093479393
$ws.Cells.Item(22, 12).Value = 1.590282990665685
$ws.Cells.Item(22, 13).Value = 1.100576827018045
$ws.Cells.Item(22, 14).Value = 2.504906168350281

$ws.Cells.Item(23, 2).Value = 3.104402380163378
$ws.Cells.Item(23, 4).Value = 0.4141278361681202
$ws.Cells.Item(23, 5).Value = 0.3009285131191746
$ws.Cells.Item(23, 6).Value = 5.986862196527198
$ws.Cells.Item(23, 7).Value = 0.002624746578415278
$ws.Cells.Item(23, 10).Value = 0.3245050554292277
$ws.Cells.Item(23, 12).Value = 1.562463177616394
$ws.Cells.Item(23, 13).Value = 1.084767401094581
$ws.Cells.Item(23, 14).Value = 2.516659709760617

$ws.Cells.Item(24, 2).Value = 2.966644089995896
$ws.Cells.Item(24, 4).Value = 0.3606193732680936
$ws.Cells.Item(24, 5).Value = 0.257360186919982
$ws.Cells.Item(24, 6).Value = 5.83686912647326
$ws.Cells.Item(24, 7).Value = 0.002639767423719719
$ws.Cells.Item(24, 10).Value = 0.2786129596527758
$ws.Cells.Item(24, 12).Value = 1.458815217370244
$ws.Cells.Item(24, 13).Value = 1.02602220863416
$ws.Cells.Item(24, 14).Value = 2.562752847321022

$ws.Cells.Item(25, 2).Value = 2.824476224194598
$ws.Cells.Item(25, 4).Value = 0.3038018865807146
$ws.Cells.Item(25, 5).Value = 0.210703920059359
$ws.Cells.Item(25, 6).Value = 5.689844739609043
$ws.Cells.Item(25, 7).Value = 0.002657107864858827
$ws.Cells.Item(25, 10).Value = 0.2292252998244209
$ws.Cells.Item(25, 12).Value = 1.350563783354318
$ws.Cells.Item(25, 13).Value = 0.9650156659979672
$ws.Cells.Item(25, 14).Value = 2.615821491174643
